$wb = $excel.ActiveWorkbook

# --- Thunderbird: fill in the new Lucene-LDA-derived keyword metric columns (C:F) ---
# for rows 45-48, 51-55, 57-61 (row 49 intentionally left without metrics, as in source data)
$ws = $wb.Worksheets.Item("Thunderbird")
$ws.Range("C45").Value = '0.431 0.343 0.480 0.480 0.384'
$ws.Range("D45").Value = '0.601 0.513 0.073 0.074 0.826'
$ws.Range("E45").Value = '0.280 0.209 0.603 0.601 0.238'
$ws.Range("F45").Value = '0.781 0.845 0.414 0.414 0.979'
$ws.Range("C46").Value = '0.673 0.713 0.655 0.710 0.574'
$ws.Range("D46").Value = '0.403 0.341 0.192 0.202 0.158'
$ws.Range("E46").Value = '0.656 0.674 0.545 0.640 0.413'
$ws.Range("F46").Value = '0.682 0.744 0.800 0.786 0.928'
$ws.Range("C47").Value = '0.588 0.583 0.326 0.408 0.475'
$ws.Range("D47").Value = '0.532 0.547 0.265 0.333 0.403 '
$ws.Range("E47").Value = '0.443 0.424 0.196 0.259 0.313 '
$ws.Range("F47").Value = '0.771 0.855 0.902 0.907 0.970'
$ws.Range("C48").Value = '0.355 0.286 0.000 0.009 0.280'
$ws.Range("D48").Value = '0.500 0.519 0.000 0.167 0.765'
$ws.Range("E48").Value = '0.220 0.167 0.000 0.004 0.163 '
$ws.Range("F48").Value = '0.758 0.845 0.926 0.925 0.977 '
$ws.Range("C51").Value = '0.431 0.410 0.423 0.424 0.416'
$ws.Range("D51").Value = '0.584 0.560 0.074 0.075 0.724 '
$ws.Range("E51").Value = '0.280 0.260 0.714 0.715 0.263'
$ws.Range("F51").Value = '0.777 0.852 0.331 0.332 0.978'
$ws.Range("C52").Value = '0.670 0.720 0.642 0.701 0.597 '
$ws.Range("D52").Value = '0.430 0.367 0.206 0.216 0.160 '
$ws.Range("E52").Value = '0.612 0.663 0.518 0.610 0.438'
$ws.Range("F52").Value = '0.710 0.768 0.820 0.808 0.926'
$ws.Range("C53").Value = '0.602 0.595 0.320 0.439 0.489'
$ws.Range("D53").Value = '0.534 0.531 0.239 0.342 0.456'
$ws.Range("E53").Value = '0.459 0.438 0.192 0.285 0.325'
$ws.Range("F53").Value = '0.772 0.851 0.897 0.907 0.972'
$ws.Range("C54").Value = '0.384 0.277 0.000 0.000 0.298 '
$ws.Range("D54").Value = '0.523 0.523 0.000 0.000 0.667 '
$ws.Range("E54").Value = '0.242 0.161 0.000 0.000 0.175 '
$ws.Range("F54").Value = '0.763 0.845 0.927 0.924 0.976'
$ws.Range("C55").Value = '0.341 0.287 0.149 0.253 0.416'
$ws.Range("D55").Value = '0.790 0.920 0.818 0.943 1.000'
$ws.Range("E55").Value = '0.206 0.167 0.080 0.145 0.263'
$ws.Range("F55").Value = '0.795 0.867 0.932 0.936 0.981 '
$ws.Range("C57").Value = '0.479 0.340 0.478 0.322 0.298'
$ws.Range("D57").Value = '0.607 0.532 0.072 0.076 0.737'
$ws.Range("E57").Value = '0.323 0.207 0.598 0.820 0.175 '
$ws.Range("F57").Value = '0.785 0.847 0.413 0.246 0.977'
$ws.Range("C58").Value = '0.673 0.703 0.572 0.644 0.611'
$ws.Range("D58").Value = '0.458 0.385 0.195 0.233 0.205'
$ws.Range("E58").Value = '0.594 0.618 0.429 0.513 0.450 '
$ws.Range("F58").Value = '0.731 0.785 0.830 0.839 0.940'
$ws.Range("C59").Value = '0.607 0.581 0.301 0.450 0.504'
$ws.Range("D59").Value = '0.555 0.537 0.244 0.372 0.519 '
$ws.Range("E59").Value = '0.463 0.421 0.179 0.294 0.338'
$ws.Range("F59").Value = '0.780 0.852 0.900 0.911 0.975'
$ws.Range("C60").Value = '0.387 0.246 0.009 0.034 0.261'
$ws.Range("D60").Value = '0.513 0.523 0.200 0.667 0.800'
$ws.Range("E60").Value = '0.245 0.140 0.004 0.018 0.150'
$ws.Range("F60").Value = '0.761 0.845 0.926 0.927 0.977'
$ws.Range("C61").Value = '0.347 0.296 0.194 0.260 0.416'
$ws.Range("D61").Value = '0.867 0.913 0.800 0.944 1.000 '
$ws.Range("E61").Value = '0.210 0.174 0.107 0.149 0.263 '
$ws.Range("F61").Value = '0.801 0.868 0.933 0.936 0.981'


# --- Cosmetic: restore the viewport/selection state captured when the
#     workbook was last saved after this edit ---
$wsLucene = $wb.Worksheets.Item("Lucene")
$wsLucene.Activate()
$wsLucene.Range("B47").Select()

$ws.Activate()
$ws.Range("C48").Select()
